$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates (subject id columns)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updated meanEMG values
$ws.Range("B2").Value = 169.82593488081545
$ws.Range("C2").Value = 136.81760523438663
$ws.Range("D2").Value = 170.77378438415633
$ws.Range("E2").Value = 134.85450354383354

# Row 3 (STR) updated meanEMG values
$ws.Range("B3").Value = 153.89250429911436
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 158.79461689412258
$ws.Range("E3").Value = 132.33390606648163

# Update selection to match the new authored range
$ws.Range("B1:E3").Select()
